$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4 (anchor G=5470)
$ws.Range("H4").Value = 600
$ws.Range("I4").Value = 600
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -486
$ws.Range("N4").ClearContents()

# Row 33 (anchor G=5512)
$ws.Range("H33").Value = 841.93335
$ws.Range("I33").Value = 844.9286
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 844.9286
$ws.Range("L33").Value = 800
$ws.Range("M33").Value = -615.9286
$ws.Range("N33").Value = -1258

# Row 76 (anchor G=12602)
$ws.Range("H76").Value = 2932.3447
$ws.Range("I76").Value = 2685.75
$ws.Range("J76").Value = 4116
$ws.Range("K76").Value = 2685.75
$ws.Range("L76").Value = 4116
$ws.Range("M76").Value = -2370.75
$ws.Range("N76").Value = -4746

# Row 79 (anchor G=12602)
$ws.Range("H79").Value = 2932.3447
$ws.Range("I79").Value = 2685.75
$ws.Range("J79").Value = 4116
$ws.Range("K79").Value = 2685.75
$ws.Range("L79").Value = 4116
$ws.Range("M79").Value = -1593.75
$ws.Range("N79").Value = -6300

# Row 86 (anchor G=12603)
$ws.Range("H86").Value = 137389.22
$ws.Range("I86").Value = 246000.6
$ws.Range("J86").Value = 1625
$ws.Range("K86").Value = 246000.6
$ws.Range("L86").Value = 1625
$ws.Range("M86").Value = -244877.6
$ws.Range("N86").Value = -3871

# Row 89 (anchor G=12603)
$ws.Range("H89").Value = 137389.22
$ws.Range("I89").Value = 246000.6
$ws.Range("J89").Value = 1625
$ws.Range("K89").Value = 1230003
$ws.Range("L89").Value = 8125
$ws.Range("M89").Value = -1224387
$ws.Range("N89").Value = -19357

# Row 106 (anchor G=19903)
$ws.Range("H106").Value = 3725.7144
$ws.Range("I106").Value = 3200
$ws.Range("J106").Value = 3936
$ws.Range("K106").Value = 3200
$ws.Range("L106").Value = 3936
$ws.Range("M106").Value = -2569
$ws.Range("N106").Value = -5198

# Row 138 (anchor G=44169)
$ws.Range("H138").Value = 2149.04
$ws.Range("I138").Value = 1807.1765
$ws.Range("J138").Value = 2875.5
$ws.Range("K138").Value = 5421.529500000001
$ws.Range("L138").Value = 8626.5
$ws.Range("M138").Value = -281.5295000000006
$ws.Range("N138").Value = -18906.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 3 (anchor G=2494)
$ws.Range("H3").Value = 401.66666
$ws.Range("I3").Value = 401.66666
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 401.66666
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -286.66666
$ws.Range("N3").ClearContents()

# Row 32 (anchor G=44147)
$ws.Range("H32").Value = 7237.763
$ws.Range("I32").Value = 6736.8887
$ws.Range("K32").Value = 6736.8887
$ws.Range("M32").Value = -6449.8887

# Row 124 (anchor G=34252)
$ws.Range("H124").Value = 21959.777
$ws.Range("J124").Value = 21959.777
$ws.Range("L124").Value = 21959.777
$ws.Range("N124").Value = -31779.777

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 7 (anchor G=1602)
$ws.Range("H7").Value = 401
$ws.Range("I7").Value = 401
$ws.Range("K7").Value = 401
$ws.Range("M7").Value = -288

# Row 86 (anchor G=12526)
$ws.Range("H86").Value = 1128.3889
$ws.Range("I86").Value = 924.8889
$ws.Range("J86").Value = 1331.8889
$ws.Range("K86").Value = 924.8889
$ws.Range("L86").Value = 1331.8889
$ws.Range("M86").Value = 198.1111
$ws.Range("N86").Value = -3577.8889

# Row 89 (anchor G=12526)
$ws.Range("H89").Value = 1128.3889
$ws.Range("I89").Value = 924.8889
$ws.Range("J89").Value = 1331.8889
$ws.Range("K89").Value = 4624.444500000001
$ws.Range("L89").Value = 6659.4445
$ws.Range("M89").Value = 991.5554999999995
$ws.Range("N89").Value = -17891.4445

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2 (anchor G=1820)
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187

# Row 31 (anchor G=44023)
$ws.Range("H31").Value = 5378155
$ws.Range("I31").Value = 1474.3489
$ws.Range("J31").Value = 17546432
$ws.Range("K31").Value = 1474.3489
$ws.Range("L31").Value = 17546432
$ws.Range("M31").Value = -1179.3489
$ws.Range("N31").Value = -17547022

# Row 34 (anchor G=44023)
$ws.Range("H34").Value = 5378155
$ws.Range("I34").Value = 1474.3489
$ws.Range("J34").Value = 17546432
$ws.Range("K34").Value = 1474.3489
$ws.Range("L34").Value = 17546432
$ws.Range("M34").Value = -1272.3489
$ws.Range("N34").Value = -17546836

# Row 62 (anchor G=12580)
$ws.Range("H62").Value = 22110.5
$ws.Range("I62").Value = 16515
$ws.Range("J62").Value = 35166.668
$ws.Range("K62").Value = 16515
$ws.Range("L62").Value = 35166.668
$ws.Range("M62").Value = -15891
$ws.Range("N62").Value = -36414.668

# Row 65 (anchor G=12580)
$ws.Range("H65").Value = 22110.5
$ws.Range("I65").Value = 16515
$ws.Range("J65").Value = 35166.668
$ws.Range("K65").Value = 82575
$ws.Range("L65").Value = 175833.34
$ws.Range("M65").Value = -79455
$ws.Range("N65").Value = -182073.34

# Row 74 (anchor G=10636)
$ws.Range("H74").Value = 43499.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 43499.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 43499.5
$ws.Range("N74").Value = -45247.5
$ws.Range("M74").ClearContents()

# Row 77 (anchor G=10636)
$ws.Range("H77").Value = 43499.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 43499.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 130498.5
$ws.Range("N77").Value = -139234.5
$ws.Range("M77").ClearContents()

# Row 88 (anchor G=10608)
$ws.Range("H88").Value = 29066.666
$ws.Range("J88").Value = 29066.666
$ws.Range("L88").Value = 29066.666
$ws.Range("N88").Value = -29878.666

# Row 91 (anchor G=10608)
$ws.Range("H91").Value = 29066.666
$ws.Range("J91").Value = 29066.666
$ws.Range("L91").Value = 29066.666
$ws.Range("N91").Value = -31874.666

# Row 92 (anchor G=18041)
$ws.Range("H92").Value = 13800.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 13800.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 13800.5
$ws.Range("N92").Value = -18792.5
$ws.Range("M92").ClearContents()

# Row 96 (anchor G=18193)
$ws.Range("H96").Value = 16100
$ws.Range("J96").Value = 16100
$ws.Range("L96").Value = 16100
$ws.Range("N96").Value = -21592

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4 (anchor G=4650)
$ws.Range("H4").Value = 3996.4285
$ws.Range("I4").Value = 12712.5
$ws.Range("J4").Value = 510
$ws.Range("K4").Value = 38137.5
$ws.Range("L4").Value = 1530
$ws.Range("M4").Value = -38025.5
$ws.Range("N4").Value = -1754

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2 (anchor G=5062)
$ws.Range("H2").Value = 57.785713
$ws.Range("I2").Value = 49.75
$ws.Range("J2").Value = 68.5
$ws.Range("K2").Value = 49.75
$ws.Range("L2").Value = 68.5
$ws.Range("M2").Value = 63.25
$ws.Range("N2").Value = -294.5

# Row 80 (anchor G=12521)
$ws.Range("H80").Value = 3231.6
$ws.Range("I80").Value = 2848.5
$ws.Range("J80").Value = 3997.8
$ws.Range("K80").Value = 2848.5
$ws.Range("L80").Value = 3997.8
$ws.Range("M80").Value = -1850.5
$ws.Range("N80").Value = -5993.8

# Row 83 (anchor G=12521)
$ws.Range("H83").Value = 3231.6
$ws.Range("I83").Value = 2848.5
$ws.Range("J83").Value = 3997.8
$ws.Range("K83").Value = 14242.5
$ws.Range("L83").Value = 19989
$ws.Range("M83").Value = -9250.5
$ws.Range("N83").Value = -29973

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16 (anchor G=5289)
$ws.Range("H16").Value = 2439.4211
$ws.Range("I16").Value = 549.94116
$ws.Range("J16").Value = 18500
$ws.Range("K16").Value = 549.94116
$ws.Range("L16").Value = 18500
$ws.Range("M16").Value = -379.94116
$ws.Range("N16").Value = -18840

# Row 46 (anchor G=5282)
$ws.Range("H46").Value = 834513.3
$ws.Range("I46").Value = 967.1429000000001
$ws.Range("J46").Value = 2001478
$ws.Range("K46").Value = 967.1429000000001
$ws.Range("L46").Value = 2001478
$ws.Range("M46").Value = -779.1429000000001
$ws.Range("N46").Value = -2001854

# Row 55 (anchor G=5284)
$ws.Range("H55").Value = 448
$ws.Range("I55").Value = 491.42856
$ws.Range("J55").Value = 404.57144
$ws.Range("K55").Value = 491.42856
$ws.Range("L55").Value = 404.57144
$ws.Range("M55").Value = -318.42856
$ws.Range("N55").Value = -750.5714399999999

# Row 68 (anchor G=12563)
$ws.Range("H68").Value = 2351.75
$ws.Range("I68").Value = 2022.1
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2022.1
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1273.1
$ws.Range("N68").Value = -5498

# Row 71 (anchor G=12563)
$ws.Range("H71").Value = 2351.75
$ws.Range("I71").Value = 2022.1
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 10110.5
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -6366.5
$ws.Range("N71").Value = -27488

# Row 82 (anchor G=12565)
$ws.Range("H82").Value = 1641.2222
$ws.Range("I82").Value = 1077.2
$ws.Range("J82").Value = 2346.25
$ws.Range("K82").Value = 1077.2
$ws.Range("L82").Value = 2346.25
$ws.Range("M82").Value = -716.2
$ws.Range("N82").Value = -3068.25

# Row 85 (anchor G=12565)
$ws.Range("H85").Value = 1641.2222
$ws.Range("I85").Value = 1077.2
$ws.Range("J85").Value = 2346.25
$ws.Range("K85").Value = 1077.2
$ws.Range("L85").Value = 2346.25
$ws.Range("M85").Value = 170.8
$ws.Range("N85").Value = -4842.25

# Row 127 (anchor G=34401)
$ws.Range("H127").Value = 22571.666
$ws.Range("J127").Value = 22571.666
$ws.Range("L127").Value = 22571.666
$ws.Range("N127").Value = -32491.666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 97 (anchor G=18220)
$ws.Range("H97").Value = 22990.666
$ws.Range("J97").Value = 22990.666
$ws.Range("L97").Value = 22990.666
$ws.Range("N97").Value = -24972.666

# Row 98 (anchor G=18374)
$ws.Range("H98").Value = 45000
$ws.Range("J98").Value = 45000
$ws.Range("L98").Value = 45000
$ws.Range("N98").Value = -50990

# Row 122 (anchor G=36208)
$ws.Range("H122").Value = 3099.1428
$ws.Range("I122").Value = 938.8
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 2816.4
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -366.3999999999996
$ws.Range("N122").Value = -30400
